$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 175.54546
$ws.Range("I33").Value = 175.54546
$ws.Range("K33").Value = 175.54546
$ws.Range("M33").Value = 53.45454000000001

$ws.Range("H53").Value = 1504.2858
$ws.Range("I53").Value = 2572
$ws.Range("J53").Value = 80.666664
$ws.Range("K53").Value = 2572
$ws.Range("L53").Value = 80.666664
$ws.Range("M53").Value = -1935
$ws.Range("N53").Value = -1354.666664

$ws.Range("H88").Value = 2342.5715
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2483
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2483
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3295

$ws.Range("H91").Value = 2342.5715
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2483
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2483
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5291

$ws.Range("H134").Value = 59780
$ws.Range("J134").Value = 59780
$ws.Range("L134").Value = 59780
$ws.Range("N134").Value = -69920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 464408.75
$ws.Range("I2").Value = 795000.5600000001
$ws.Range("J2").Value = 1580.2
$ws.Range("K2").Value = 795000.5600000001
$ws.Range("L2").Value = 1580.2
$ws.Range("M2").Value = -794887.5600000001
$ws.Range("N2").Value = -1806.2

$ws.Range("H45").Value = 1347.9412
$ws.Range("J45").Value = 1735.2
$ws.Range("L45").Value = 1735.2
$ws.Range("N45").Value = -2489.2

$ws.Range("H61").Value = 9842.462
$ws.Range("I61").Value = 14859.333
$ws.Range("K61").Value = 14859.333
$ws.Range("M61").Value = -14647.333

$ws.Range("H74").Value = 1223.84
$ws.Range("I74").Value = 494.95
$ws.Range("K74").Value = 494.95
$ws.Range("M74").Value = 379.05

$ws.Range("H77").Value = 1223.84
$ws.Range("I77").Value = 494.95
$ws.Range("K77").Value = 2474.75
$ws.Range("M77").Value = 1893.25

$ws.Range("H116").Value = 464408.75
$ws.Range("I116").Value = 795000.5600000001
$ws.Range("J116").Value = 1580.2
$ws.Range("K116").Value = 795000.5600000001
$ws.Range("L116").Value = 1580.2
$ws.Range("M116").Value = -792706.5600000001
$ws.Range("N116").Value = -6168.2

$ws.Range("H136").Value = 9842.462
$ws.Range("I136").Value = 14859.333
$ws.Range("K136").Value = 44577.999
$ws.Range("M136").Value = -42027.999

$ws.Range("H139").Value = 40920
$ws.Range("J139").Value = 40920
$ws.Range("L139").Value = 40920
$ws.Range("N139").Value = -51200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 464408.75
$ws.Range("I3").Value = 795000.5600000001
$ws.Range("J3").Value = 1580.2
$ws.Range("K3").Value = 795000.5600000001
$ws.Range("L3").Value = 1580.2
$ws.Range("M3").Value = -794886.5600000001
$ws.Range("N3").Value = -1808.2

$ws.Range("H86").Value = 186250.9
$ws.Range("I86").Value = 6458.3335
$ws.Range("J86").Value = 402002
$ws.Range("K86").Value = 6458.3335
$ws.Range("L86").Value = 402002
$ws.Range("M86").Value = -5335.3335
$ws.Range("N86").Value = -404248

$ws.Range("H89").Value = 186250.9
$ws.Range("I89").Value = 6458.3335
$ws.Range("J89").Value = 402002
$ws.Range("K89").Value = 32291.6675
$ws.Range("L89").Value = 2010010
$ws.Range("M89").Value = -26675.6675
$ws.Range("N89").Value = -2021242

$ws.Range("H105").Value = 2342.2693
$ws.Range("I105").Value = 2039.0869
$ws.Range("K105").Value = 2039.0869
$ws.Range("M105").Value = -292.0869

$ws.Range("H107").Value = 896
$ws.Range("I107").Value = 679.5454999999999
$ws.Range("K107").Value = 679.5454999999999
$ws.Range("M107").Value = 1240.4545

$ws.Range("H135").Value = 57999.5
$ws.Range("J135").Value = 57999.5
$ws.Range("L135").Value = 57999.5
$ws.Range("N135").Value = -68139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30782

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30364

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H132").Value = 2590.1428
$ws.Range("I132").Value = 1692.3334
$ws.Range("K132").Value = 5077.0002
$ws.Range("M132").Value = -2547.0002

$ws.Range("H134").Value = 2677.05
$ws.Range("I134").Value = 2336.1765
$ws.Range("K134").Value = 7008.529500000001
$ws.Range("M134").Value = -4473.529500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13842.2295
$ws.Range("J131").Value = 15021.714
$ws.Range("L131").Value = 45065.142
$ws.Range("N131").Value = -55145.142

$ws.Range("H132").Value = 1633.7273
$ws.Range("J132").Value = 1749.4445
$ws.Range("L132").Value = 15745.0005
$ws.Range("N132").Value = -20805.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 49900
$ws.Range("J134").Value = 49900
$ws.Range("L134").Value = 49900
$ws.Range("N134").Value = -60040

$ws.Range("H136").Value = 3317.5
$ws.Range("I136").Value = 3500.25
$ws.Range("J136").Value = 2952
$ws.Range("K136").Value = 10500.75
$ws.Range("L136").Value = 8856
$ws.Range("M136").Value = -7950.75
$ws.Range("N136").Value = -13956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1334.7678
$ws.Range("I132").Value = 968.63416
$ws.Range("J132").Value = 2335.5334
$ws.Range("K132").Value = 2905.90248
$ws.Range("L132").Value = 7006.600199999999
$ws.Range("M132").Value = -375.9024799999997
$ws.Range("N132").Value = -12066.6002

$ws.Range("H133").Value = 65486.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 65486.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 65486.25
$ws.Range("N133").Value = -75606.25
$ws.Range("M133").ClearContents()
